$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3240.4
$ws.Range("I32").Value = 3543
$ws.Range("J32").Value = 2786.5
$ws.Range("K32").Value = 3543
$ws.Range("L32").Value = 2786.5
$ws.Range("M32").Value = -3217
$ws.Range("N32").Value = -3438.5
$ws.Range("H86").Value = 46066.9
$ws.Range("I86").Value = 1867.25
$ws.Range("J86").Value = 75533.336
$ws.Range("K86").Value = 1867.25
$ws.Range("L86").Value = 75533.336
$ws.Range("M86").Value = -744.25
$ws.Range("N86").Value = -77779.336
$ws.Range("H89").Value = 46066.9
$ws.Range("I89").Value = 1867.25
$ws.Range("J89").Value = 75533.336
$ws.Range("K89").Value = 9336.25
$ws.Range("L89").Value = 377666.68
$ws.Range("M89").Value = -3720.25
$ws.Range("N89").Value = -388898.68
$ws.Range("H94").Value = 50129930
$ws.Range("I94").Value = 83344050
$ws.Range("J94").Value = 308751.5
$ws.Range("K94").Value = 83344050
$ws.Range("L94").Value = 308751.5
$ws.Range("M94").Value = -83343599
$ws.Range("N94").Value = -309653.5
$ws.Range("H116").Value = 20376352
$ws.Range("J116").Value = 8971.75
$ws.Range("L116").Value = 8971.75
$ws.Range("N116").Value = -15855.75
$ws.Range("H132").Value = 2882.8542
$ws.Range("I132").Value = 2577.8262
$ws.Range("K132").Value = 7733.4786
$ws.Range("M132").Value = -5203.4786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8150.1763
$ws.Range("I2").Value = 9599.23
$ws.Range("K2").Value = 9599.23
$ws.Range("M2").Value = -9486.23
$ws.Range("H63").Value = 2749.75
$ws.Range("I63").Value = 1999.6666
$ws.Range("K63").Value = 1999.6666
$ws.Range("M63").Value = -1313.6666
$ws.Range("H66").Value = 2749.75
$ws.Range("I66").Value = 1999.6666
$ws.Range("K66").Value = 9998.333000000001
$ws.Range("M66").Value = -6566.333000000001
$ws.Range("H116").Value = 8150.1763
$ws.Range("I116").Value = 9599.23
$ws.Range("K116").Value = 9599.23
$ws.Range("M116").Value = -7305.23
$ws.Range("H132").Value = 3593.258
$ws.Range("I132").Value = 2449.0527
$ws.Range("K132").Value = 7347.158100000001
$ws.Range("M132").Value = -4817.158100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8150.1763
$ws.Range("I3").Value = 9599.23
$ws.Range("K3").Value = 9599.23
$ws.Range("M3").Value = -9485.23
$ws.Range("H105").Value = 2936.1875
$ws.Range("I105").Value = 1880.75
$ws.Range("J105").Value = 6102.5
$ws.Range("K105").Value = 1880.75
$ws.Range("L105").Value = 6102.5
$ws.Range("M105").Value = -133.75
$ws.Range("N105").Value = -9596.5
$ws.Range("H107").Value = 4356.6294
$ws.Range("I107").Value = 4892.227
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 4892.227
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -2972.227
$ws.Range("N107").Value = -5840
$ws.Range("H134").Value = 2244.1587
$ws.Range("I134").Value = 2079.7273
$ws.Range("K134").Value = 6239.1819
$ws.Range("M134").Value = -3704.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1712.375
$ws.Range("I16").Value = 1549.8334
$ws.Range("J16").Value = 2200
$ws.Range("K16").Value = 1549.8334
$ws.Range("L16").Value = 2200
$ws.Range("M16").Value = -1262.8334
$ws.Range("N16").Value = -2774
$ws.Range("H95").Value = 312031200
$ws.Range("J95").Value = 312031200
$ws.Range("L95").Value = 312031200
$ws.Range("N95").Value = -312036692
$ws.Range("H99").Value = 7823369.5
$ws.Range("J99").Value = 4544.2856
$ws.Range("L99").Value = 4544.2856
$ws.Range("N99").Value = -7540.2856
$ws.Range("H105").Value = 16181.286
$ws.Range("I105").Value = 18044.834
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 18044.834
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -16297.834
$ws.Range("N105").Value = -8494
$ws.Range("H107").Value = 7425.6113
$ws.Range("I107").Value = 9523.885
$ws.Range("K107").Value = 9523.885
$ws.Range("M107").Value = -7603.885
$ws.Range("H113").Value = 1712.375
$ws.Range("I113").Value = 1549.8334
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1549.8334
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 620.1666
$ws.Range("N113").Value = -6540
$ws.Range("H126").Value = 7823369.5
$ws.Range("J126").Value = 4544.2856
$ws.Range("L126").Value = 13632.8568
$ws.Range("N126").Value = -18572.8568
$ws.Range("H132").Value = 18960.297
$ws.Range("I132").Value = 2477.12
$ws.Range("J132").Value = 225000
$ws.Range("K132").Value = 7431.36
$ws.Range("L132").Value = 675000
$ws.Range("M132").Value = -4901.36
$ws.Range("N132").Value = -680060
$ws.Range("H134").Value = 4171466
$ws.Range("I134").Value = 5212749.5
$ws.Range("K134").Value = 15638248.5
$ws.Range("M134").Value = -15635713.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 950
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 950
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 2850
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -3426
$ws.Range("H121").Value = 4665762.5
$ws.Range("J121").Value = 6666861
$ws.Range("L121").Value = 20000583
$ws.Range("N121").Value = -20003203

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4850.5454
$ws.Range("I70").Value = 4498.7144
$ws.Range("J70").Value = 5466.25
$ws.Range("K70").Value = 4498.7144
$ws.Range("L70").Value = 5466.25
$ws.Range("M70").Value = -4228.7144
$ws.Range("N70").Value = -6006.25
$ws.Range("H73").Value = 4850.5454
$ws.Range("I73").Value = 4498.7144
$ws.Range("J73").Value = 5466.25
$ws.Range("K73").Value = 4498.7144
$ws.Range("L73").Value = 5466.25
$ws.Range("M73").Value = -3562.7144
$ws.Range("N73").Value = -7338.25
$ws.Range("H80").Value = 10750.5
$ws.Range("J80").Value = 3333
$ws.Range("L80").Value = 3333
$ws.Range("N80").Value = -5329
$ws.Range("H83").Value = 10750.5
$ws.Range("J83").Value = 3333
$ws.Range("L83").Value = 16665
$ws.Range("N83").Value = -26649
$ws.Range("H107").Value = 614.5833
$ws.Range("I107").Value = 520
$ws.Range("J107").Value = 898.3333
$ws.Range("K107").Value = 520
$ws.Range("L107").Value = 898.3333
$ws.Range("M107").Value = 1400
$ws.Range("N107").Value = -4738.3333
$ws.Range("H122").Value = 20180.375
$ws.Range("I122").Value = 17106.615
$ws.Range("K122").Value = 51319.845
$ws.Range("M122").Value = -48869.845
$ws.Range("H135").Value = 82199.7
$ws.Range("J135").Value = 82199.7
$ws.Range("L135").Value = 82199.7
$ws.Range("N135").Value = -92339.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 32033.6
$ws.Range("I7").Value = 47889.89
$ws.Range("K7").Value = 47889.89
$ws.Range("M7").Value = -47777.89
$ws.Range("H61").Value = 4170.75
$ws.Range("I61").Value = 3021.3215
$ws.Range("J61").Value = 8193.75
$ws.Range("K61").Value = 3021.3215
$ws.Range("L61").Value = 8193.75
$ws.Range("M61").Value = -2819.3215
$ws.Range("N61").Value = -8597.75
$ws.Range("H101").Value = 32706.75
$ws.Range("J101").Value = 32706.75
$ws.Range("L101").Value = 32706.75
$ws.Range("N101").Value = -39196.75
$ws.Range("H113").Value = 4170.75
$ws.Range("I113").Value = 3021.3215
$ws.Range("J113").Value = 8193.75
$ws.Range("K113").Value = 3021.3215
$ws.Range("L113").Value = 8193.75
$ws.Range("M113").Value = -851.3215
$ws.Range("N113").Value = -12533.75
$ws.Range("H126").Value = 32033.6
$ws.Range("I126").Value = 47889.89
$ws.Range("K126").Value = 143669.67
$ws.Range("M126").Value = -141199.67
$ws.Range("H136").Value = 4326.839
$ws.Range("I136").Value = 1728.3334
$ws.Range("J136").Value = 7924.769
$ws.Range("K136").Value = 5185.0002
$ws.Range("L136").Value = 23774.307
$ws.Range("M136").Value = -2635.0002
$ws.Range("N136").Value = -28874.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14714
$ws.Range("J81").Value = 2811.4285
$ws.Range("L81").Value = 5622.857
$ws.Range("N81").Value = -7744.857
$ws.Range("H84").Value = 14714
$ws.Range("J84").Value = 2811.4285
$ws.Range("L84").Value = 28114.285
$ws.Range("N84").Value = -38722.285
$ws.Range("H97").Value = 53439.4
$ws.Range("J97").Value = 53439.4
$ws.Range("L97").Value = 53439.4
$ws.Range("N97").Value = -55421.4
$ws.Range("H107").Value = 16093.904
$ws.Range("I107").Value = 1967.7858
$ws.Range("J107").Value = 44346.145
$ws.Range("K107").Value = 5903.357400000001
$ws.Range("L107").Value = 133038.435
$ws.Range("M107").Value = -3983.357400000001
$ws.Range("N107").Value = -136878.435
$ws.Range("H132").Value = 12321.484
$ws.Range("I132").Value = 14185.741
$ws.Range("K132").Value = 42557.223
$ws.Range("M132").Value = -40027.223
